$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96, pushing the existing rows 96-99 down to 97-100.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new weekly price record.
$ws.Cells.Item(96, 1).Value = 11
$ws.Cells.Item(96, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(96, 3).Value = "Bíobío"
$ws.Cells.Item(96, 4).Value = 44610
$ws.Cells.Item(96, 5).Value = 8
$ws.Cells.Item(96, 6).Value = 100112032
$ws.Cells.Item(96, 7).Value = "Zapallo italiano"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 1050
$ws.Cells.Item(96, 11).Value = 9000
$ws.Cells.Item(96, 12).Value = 10000
$ws.Cells.Item(96, 13).Value = 9048
$ws.Cells.Item(96, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(96, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(96, 16).Value = 181
$ws.Cells.Item(96, 17).Value = 50
$ws.Cells.Item(96, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the neighbouring rows.
$ws.Cells.Item(96, 4).NumberFormat = $ws.Cells.Item(97, 4).NumberFormat
